# Knowledge Distillation deck update:
#  - Slide 2 ("Topics") bullet list replaced: Teacher Model/Student Model/Model
#    Compression -> Process/Training/Algorithms/Applications/Conclusions
#  - 11 new slides inserted before the closing "Questions?" slide, covering
#    Feature Based Knowledge, Relation Based Knowledge, Knowledge Distillation
#    Training, Offline/Online/Self Distillation, Knowledge Distillation
#    Algorithms, Applications and Conclusions.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Update slide 2 ("Topics") content placeholder text
# ---------------------------------------------------------------------------
$topics = $p.Slides.Item(2)
$topicsBody = $topics.Shapes.Item(2).TextFrame.TextRange
$topicsBody.Text = "Knowledge Distillation`rProcess`rTraining`rAlgorithms`rApplications`r"

# ---------------------------------------------------------------------------
# Helper data: the closing slide currently sits at index 11 ("Questions?").
# We insert the new slides ahead of it, one at a time, so it naturally ends
# up at index 22 once all eleven inserts are done.
# ---------------------------------------------------------------------------
$insertAt = 11

# --- Slide 11: "Feature Based Knowledge" (title only) ----------------------
$s = $p.Slides.Add($insertAt, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Feature Based Knowledge"
$s.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2
$s.Shapes.Item(2).Delete()
$insertAt = $insertAt + 1

# --- Slide 12: "Feature Based Knowledge" (title + empty content) -----------
$s = $p.Slides.Add($insertAt, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Feature Based Knowledge"
$s.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2
$insertAt = $insertAt + 1

# --- Slide 13: "Relation Based Knowledge" (title only) ---------------------
$s = $p.Slides.Add($insertAt, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Relation Based Knowledge"
$s.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2
$s.Shapes.Item(2).Delete()
$insertAt = $insertAt + 1

# --- Slide 14: "Relation Based Knowledge" (title + empty content) ----------
$s = $p.Slides.Add($insertAt, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Relation Based Knowledge"
$s.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2
$insertAt = $insertAt + 1

# --- Slide 15: "Knowledge Distillation Training" ----------------------------
$s = $p.Slides.Add($insertAt, 2)
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "Knowledge Distillation "
$titleTr.InsertAfter("Training") | Out-Null
$s.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2

$bodyTr = $s.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Text = "Offline distillation`rTeacher – Pre-trained `rStudent – Trained `rOnline distillation`rTeacher – Trained `rStudent – Trained`rSelf distillation`r"
$bodyTr.Paragraphs(2,1).IndentLevel = 2
$bodyTr.Paragraphs(3,1).IndentLevel = 2
$bodyTr.Paragraphs(5,1).IndentLevel = 2
$bodyTr.Paragraphs(6,1).IndentLevel = 2
$insertAt = $insertAt + 1

# --- Slide 16: "Offline Distillation" (title + empty content) --------------
$s = $p.Slides.Add($insertAt, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Offline Distillation"
$s.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2
$insertAt = $insertAt + 1

# --- Slide 17: "Online Distillation" (title + empty content) ---------------
$s = $p.Slides.Add($insertAt, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Online Distillation"
$s.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2
$insertAt = $insertAt + 1

# --- Slide 18: "Self Distillation" (title + empty content) -----------------
$s = $p.Slides.Add($insertAt, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Self Distillation"
$s.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2
$insertAt = $insertAt + 1

# --- Slide 19: "Knowledge Distillation Algorithms" --------------------------
$s = $p.Slides.Add($insertAt, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Knowledge Distillation Algorithms"
$s.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Adversarial Distillation`rMulti-Teacher Distillation`rCross-modal Distillation"
$insertAt = $insertAt + 1

# --- Slide 20: "Applications" (title + empty content) -----------------------
$s = $p.Slides.Add($insertAt, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Applications"
$s.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2
$insertAt = $insertAt + 1

# --- Slide 21: "Conclusions" (title + empty content) ------------------------
$s = $p.Slides.Add($insertAt, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusions"
$s.Shapes.Item(1).TextFrame.TextRange.ParagraphFormat.Alignment = 2
$insertAt = $insertAt + 1

Write-Output ("Final slide count = " + $p.Slides.Count)
